$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
$ws.Range("H11").Value = 48.57143
$ws.Range("I11").Value = 48.57143
$ws.Range("K11").Value = 48.57143
$ws.Range("M11").Value = 91.42857000000001
$ws.Range("H76").Value = 3612.9375
$ws.Range("I76").Value = 3480.3
$ws.Range("J76").Value = 3834
$ws.Range("K76").Value = 3480.3
$ws.Range("L76").Value = 3834
$ws.Range("M76").Value = -3165.3
$ws.Range("N76").Value = -4464
$ws.Range("H79").Value = 3612.9375
$ws.Range("I79").Value = 3480.3
$ws.Range("J79").Value = 3834
$ws.Range("K79").Value = 3480.3
$ws.Range("L79").Value = 3834
$ws.Range("M79").Value = -2388.3
$ws.Range("N79").Value = -6018
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null
$ws.Range("H112").Value = 3781.818
$ws.Range("J112").Value = 4050
$ws.Range("L112").Value = 12150
$ws.Range("N112").Value = -14366
$ws.Range("H113").Value = 3875
$ws.Range("J113").Value = 4000
$ws.Range("L113").Value = 4000
$ws.Range("N113").Value = -10508
$ws.Range("H129").Value = 848.98505
$ws.Range("I129").Value = 449.42856
$ws.Range("J129").Value = 895.6
$ws.Range("K129").Value = 1348.28568
$ws.Range("L129").Value = 2686.8
$ws.Range("M129").Value = 3651.71432
$ws.Range("N129").Value = -12686.8
$ws.Range("H141").Value = 862.5
$ws.Range("I141").Value = 862.5
$ws.Range("K141").Value = 2587.5
$ws.Range("M141").Value = 2592.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4474.5684
$ws.Range("J32").Value = 5456
$ws.Range("L32").Value = 5456
$ws.Range("N32").Value = -6030
$ws.Range("H94").Value = 21591
$ws.Range("J94").Value = 21591
$ws.Range("L94").Value = 21591
$ws.Range("N94").Value = -23393
$ws.Range("H102").Value = 23811538
$ws.Range("I102").Value = 27779296
$ws.Range("K102").Value = 27779296
$ws.Range("M102").Value = -27777674
$ws.Range("H104").Value = 58484.5
$ws.Range("J104").Value = 58484.5
$ws.Range("L104").Value = 58484.5
$ws.Range("N104").Value = -65472.5
$ws.Range("H110").Value = 1675.7273
$ws.Range("I110").Value = 874.2857
$ws.Range("K110").Value = 874.2857
$ws.Range("M110").Value = 1170.7143
$ws.Range("H122").Value = 1865.091
$ws.Range("I122").Value = 1851.6
$ws.Range("K122").Value = 5554.799999999999
$ws.Range("M122").Value = -3104.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2845.359
$ws.Range("I86").Value = 3031.6206
$ws.Range("J86").Value = 2305.2
$ws.Range("K86").Value = 3031.6206
$ws.Range("L86").Value = 2305.2
$ws.Range("M86").Value = -1908.6206
$ws.Range("N86").Value = -4551.2
$ws.Range("H89").Value = 2845.359
$ws.Range("I89").Value = 3031.6206
$ws.Range("J89").Value = 2305.2
$ws.Range("K89").Value = 15158.103
$ws.Range("L89").Value = 11526
$ws.Range("M89").Value = -9542.103000000001
$ws.Range("N89").Value = -22758

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2099.6086
$ws.Range("J31").Value = 2999.75
$ws.Range("L31").Value = 2999.75
$ws.Range("N31").Value = -3589.75
$ws.Range("H34").Value = 2099.6086
$ws.Range("J34").Value = 2999.75
$ws.Range("L34").Value = 2999.75
$ws.Range("N34").Value = -3403.75
$ws.Range("H43").Value = 24018.666
$ws.Range("J43").Value = 24018.666
$ws.Range("L43").Value = 24018.666
$ws.Range("N43").Value = -24386.666
$ws.Range("H101").Value = 24018.666
$ws.Range("J101").Value = 24018.666
$ws.Range("L101").Value = 24018.666
$ws.Range("N101").Value = -30508.666
$ws.Range("H107").Value = 557.5172
$ws.Range("I107").Value = 443.73685
$ws.Range("J107").Value = 773.7
$ws.Range("K107").Value = 443.73685
$ws.Range("L107").Value = 773.7
$ws.Range("M107").Value = 1476.26315
$ws.Range("N107").Value = -4613.7
$ws.Range("H132").Value = 9620.143
$ws.Range("I132").Value = 11678.4
$ws.Range("K132").Value = 35035.2
$ws.Range("M132").Value = -32505.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16130303
$ws.Range("J131").Value = 1411.6111
$ws.Range("L131").Value = 4234.8333
$ws.Range("N131").Value = -14314.8333
$ws.Range("H133").Value = 2314.125
$ws.Range("I133").Value = 928
$ws.Range("J133").Value = 2944.182
$ws.Range("K133").Value = 2784
$ws.Range("L133").Value = 8832.545999999998
$ws.Range("M133").Value = 2276
$ws.Range("N133").Value = -18952.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 89689.7
$ws.Range("I122").Value = 1700.909
$ws.Range("J122").Value = 251002.5
$ws.Range("K122").Value = 5102.727000000001
$ws.Range("L122").Value = 753007.5
$ws.Range("M122").Value = -2652.727000000001
$ws.Range("N122").Value = -757907.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 6000
$ws.Range("I58").Value = 6000
$ws.Range("K58").Value = 6000
$ws.Range("M58").Value = -5740
$ws.Range("H93").Value = 528.7143
$ws.Range("I93").Value = 528.7143
$ws.Range("K93").Value = 528.7143
$ws.Range("M93").Value = 719.2857
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null
$ws.Range("H130").Value = 69596
$ws.Range("J130").Value = 69596
$ws.Range("L130").Value = 69596
$ws.Range("N130").Value = -79636

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 10040.5
$ws.Range("J44").Value = 10040.5
$ws.Range("L44").Value = 10040.5
$ws.Range("N44").Value = -11148.5
$ws.Range("H51").Value = 9633
$ws.Range("I51").Value = 9999
$ws.Range("J51").Value = 9450
$ws.Range("K51").Value = 9999
$ws.Range("L51").Value = 9450
$ws.Range("M51").Value = -9489
$ws.Range("N51").Value = -10470
$ws.Range("H104").Value = 20369
$ws.Range("J104").Value = 20369
$ws.Range("L104").Value = 20369
$ws.Range("N104").Value = -27357
